# Append the new experiment-run row (row 25) to the "training" log sheet,
# in advance of running experiment 2025-09-16/0002.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("training")

$row = 25

$ws.Range("A$row").Value = "2025-09-16 14:13:04"
$ws.Range("B$row").Value = "training"
$ws.Range("C$row").Value = "configs/training/2025-09-16/x/0000"
$ws.Range("E$row").Value = "['cross_entropy', 'spectral_entropy']"
$ws.Range("F$row").Value = "[1.0, 0.01]"
$ws.Range("G$row").Value = "['torch.optim.adamw.AdamW', 'torch.optim.adam.Adam']"
$ws.Range("H$row").Value = "[0.001, 0.001]"
$ws.Range("I$row").Value = 128
$ws.Range("J$row").Value = 128
$ws.Range("K$row").Value = "general_utils.ml.training.NoImprovementStopping"
$ws.Range("L$row").Value = 8
$ws.Range("M$row").Value = 0.00001
$ws.Range("N$row").Value = 500
